$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "H2" = 2.7
    "I2" = 4.05
    "K2" = 1.82
    "L2" = 4.6
    "M2" = 1.42
    "N2" = 2.47
    "O2" = 2.2
    "P2" = 1.52
    "Q2" = 3.65
    "R2" = 1.2
    "S2" = 1.55
    "T2" = 2.15
    "U2" = 1.83
    "V2" = 1.78
    "W2" = 6.1
    "X2" = 9.5
    "Z2" = 21
    "AA2" = 19
    "AB2" = 32
    "AC2" = 6.6
    "AE2" = 14.5
    "AF2" = 80
    "AG2" = 700
    "AH2" = 9.5
    "AK2" = 75
    "AL2" = 45
    "AM2" = 50
    "I3" = 4.3
    "J3" = 2.27
    "K3" = 2.18
    "L3" = 4.5
    "W3" = 6.3
    "AB3" = 30
    "AE3" = 16.5
    "AG3" = 700
    "AH3" = 11.25
    "AI3" = 24
    "G4" = 1.22
    "H4" = 5.4
    "I4" = 11.75
    "L4" = 9
    "O4" = 1.52
    "U4" = 2.07
    "V4" = 1.6
    "W4" = 7.4
    "X4" = 6
    "Y4" = 9.25
    "Z4" = 7
    "AB4" = 32
    "AD4" = 11.5
    "AE4" = 27
    "AF4" = 150
    "AH4" = 30
    "AI4" = 100
    "AJ4" = 37
    "AK4" = 450
    "AL4" = 175
    "G5" = 2.62
    "H5" = 2.92
    "I5" = 2.75
    "J5" = 3.3
    "K5" = 1.91
    "L5" = 3.4
    "M5" = 1.5
    "N5" = 2.27
    "O5" = 2.4
    "P5" = 1.44
    "Q5" = 4.1
    "R5" = 1.15
    "S5" = 1.5
    "T5" = 2.25
    "U5" = 2.02
    "V5" = 1.62
    "W5" = 6.4
    "X5" = 11.5
    "Y5" = 10.5
    "Z5" = 30
    "AA5" = 27
    "AC5" = 6.4
    "AE5" = 18
    "AF5" = 120
    "AH5" = 6.6
    "AI5" = 12.5
    "AJ5" = 10.75
    "AK5" = 32
    "AM5" = 50
    "G6" = 2.32
    "H6" = 3.3
    "K6" = 2.12
    "R6" = 1.25
    "T6" = 2.6
    "X6" = 11
    "AI6" = 13.5
    "G7" = 2.05
    "H7" = 3.5
    "I7" = 3.2
    "J7" = 2.65
    "K7" = 2.12
    "L7" = 3.7
    "O7" = 1.82
    "Q7" = 2.9
    "T7" = 2.57
    "U7" = 1.7
    "V7" = 1.91
    "W7" = 7.7
    "X7" = 10
    "Y7" = 8.75
    "Z7" = 18.5
    "AA7" = 16
    "AB7" = 27
    "AD7" = 6.7
    "AH7" = 10
    "AI7" = 17
    "AJ7" = 11.5
    "AK7" = 40
    "AL7" = 28
    "AM7" = 35
    "G9" = 2.47
    "H9" = 3.4
    "J9" = 3
    "K9" = 2.15
    "L9" = 3.1
    "M9" = 1.24
    "N9" = 3.3
    "O9" = 1.72
    "P9" = 1.9
    "Q9" = 2.67
    "R9" = 1.36
    "U9" = 1.6
    "V9" = 2.07
    "W9" = 9.75
    "Z9" = 27
    "AA9" = 19
    "AB9" = 26
    "AC9" = 11.5
    "AD9" = 6.7
    "AE9" = 12.5
    "AF9" = 50
    "AG9" = 350
    "AH9" = 9.5
    "AJ9" = 9.75
    "AL9" = 20
    "AM9" = 27
    "G10" = 2.38
    "I10" = 3
    "J10" = 3.1
    "L10" = 3.6
    "S10" = 1.4
    "T10" = 2.75
    "U10" = 1.73
    "V10" = 2
    "X10" = 12
    "Z10" = 23
    "AC10" = 9.5
    "AF10" = 41
    "AH10" = 9.5
    "AJ10" = 11
    "AK10" = 29
    "AL10" = 23
    "G11" = 2.25
    "H11" = 2.88
    "I11" = 3.4
    "J11" = 3.1
    "K11" = 1.83
    "L11" = 4
    "M11" = 1.53
    "N11" = 2.38
    "O11" = 2.63
    "P11" = 1.44
    "Q11" = 5
    "R11" = 1.14
    "S11" = 1.62
    "T11" = 2.2
    "U11" = 2.25
    "V11" = 1.57
    "W11" = 6
    "AA11" = 23
    "AC11" = 6
    "AL11" = 34
    "AN11" = 1.11
    "AO11" = 6
    "AP11" = 2.03
    "AQ11" = 1.83
    "I12" = 2.25
    "J12" = 4.33
    "L12" = 3.1
    "Z12" = 41
    "AH12" = 6
    "AI12" = 9.5
    "G13" = 2.65
    "H13" = 3.2
    "I13" = 2.45
    "J13" = 3.25
    "K13" = 2.1
    "L13" = 3.1
    "M13" = 1.28
    "N13" = 3.4
    "O13" = 1.82
    "P13" = 1.88
    "Q13" = 2.92
    "R13" = 1.35
    "S13" = 1.4
    "T13" = 2.72
    "U13" = 1.65
    "V13" = 2.12
    "W13" = 9.5
    "X13" = 14.5
    "Y13" = 9.75
    "Z13" = 32
    "AA13" = 21
    "AB13" = 27
    "AC13" = 7.5
    "AD13" = 6.3
    "AE13" = 12.5
    "AH13" = 8.75
    "AI13" = 13
    "AJ13" = 9.25
    "AK13" = 27
    "AL13" = 19.5
    "AM13" = 27
    "AN13" = 1.06
    "AO13" = 7.5
    "G15" = 2.2
    "I15" = 2.8
    "J15" = 2.82
    "L15" = 3.4
    "X15" = 11.25
    "Z15" = 22
    "AA15" = 17.5
    "AC15" = 7.7
    "AK15" = 32
    "AO15" = 7.7
    "I16" = 3.15
    "M16" = 1.52
    "N16" = 2.37
    "X16" = 10.25
    "AA16" = 24
    "AE16" = 17
    "AH16" = 7.3
    "AM16" = 50
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
